$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking cells to remain stored as text (matches the
# source data which is written as inline strings, e.g. "300.27", "2.58%").
$textCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "D4",
    "E4",
    "D5",
    "E5",
    "D6",
    "E6",
    "D7",
    "E7",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "D11",
    "E11",
    "D12",
    "E12",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "D17",
    "E17",
    "D18",
    "E18",
    "D19",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "D22",
    "E22",
    "D23",
    "E23",
    "D24",
    "E24",
    "D25",
    "E25",
    "D38",
    "E38",
    "E39",
    "D40",
    "E40",
    "D41",
    "E41",
    "D42",
    "E42",
    "D43",
    "E43",
    "D44",
    "E44",
    "D45",
    "E45",
    "D46",
    "E46",
    "D47",
    "E47",
    "D48",
    "E48",
    "D50",
    "E50",
    "D51",
    "E51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data scraped on 2023-02-15.
$ws.Range("D2").Value = "300.27"
$ws.Range("E2").Value = "2.58%"
$ws.Range("D3").Value = "42.35"
$ws.Range("E3").Value = "4.61%"
$ws.Range("D4").Value = "5.009"
$ws.Range("E4").Value = "0.15%"
$ws.Range("D5").Value = "0.07608"
$ws.Range("E5").Value = "3.04%"
$ws.Range("D6").Value = "1.608"
$ws.Range("E6").Value = "2.24%"
$ws.Range("D7").Value = "0.9651"
$ws.Range("E7").Value = "4.49%"
$ws.Range("E8").Value = "0.26%"
$ws.Range("D9").Value = "0.1199"
$ws.Range("E9").Value = "0.84%"
$ws.Range("D10").Value = "0.1838"
$ws.Range("E10").Value = "1.52%"
$ws.Range("D11").Value = "0.09154"
$ws.Range("E11").Value = "3.92%"
$ws.Range("D12").Value = "0.04190"
$ws.Range("E12").Value = "-4.77%"
$ws.Range("E13").Value = "-0.51%"
$ws.Range("D14").Value = "0.001260"
$ws.Range("E14").Value = "-0.62%"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "0.04123"
$ws.Range("E15").Value = "5.23%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.005880"
$ws.Range("E16").Value = "1.22%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.352"
$ws.Range("E17").Value = "0.30%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "4.382"
$ws.Range("E18").Value = "2.05%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "0.3341"
$ws.Range("E19").Value = "0.73%"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "8.341"
$ws.Range("E20").Value = "5.08%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "0.1406"
$ws.Range("E21").Value = "1.20%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "0.3293"
$ws.Range("E22").Value = "11.26%"
$ws.Range("D23").Value = "0.001258"
$ws.Range("E23").Value = "-0.32%"
$ws.Range("D24").Value = "0.003931"
$ws.Range("E24").Value = "3.99%"
$ws.Range("D25").Value = "0.0001343"
$ws.Range("E25").Value = "7.31%"
$ws.Range("D38").Value = "0.02425"
$ws.Range("E38").Value = "4.04%"
$ws.Range("E39").Value = "3.89%"
$ws.Range("D40").Value = "0.006477"
$ws.Range("E40").Value = "10.99%"
$ws.Range("D41").Value = "0.007734"
$ws.Range("E41").Value = "-0.93%"
$ws.Range("D42").Value = "0.1337"
$ws.Range("E42").Value = "3.64%"
$ws.Range("D43").Value = "0.007355"
$ws.Range("E43").Value = "-0.48%"
$ws.Range("D44").Value = "0.007271"
$ws.Range("E44").Value = "-9.47%"
$ws.Range("D45").Value = "0.3013"
$ws.Range("E45").Value = "3.41%"
$ws.Range("D46").Value = "0.00006314"
$ws.Range("E46").Value = "1.71%"
$ws.Range("D47").Value = "0.00000000746"
$ws.Range("E47").Value = "-0.61%"
$ws.Range("D48").Value = "0.08533"
$ws.Range("E48").Value = "84.39%"
$ws.Range("D50").Value = "0.00002089"
$ws.Range("E50").Value = "-0.61%"
$ws.Range("D51").Value = "0.0001990"
$ws.Range("E51").Value = "-0.61%"
